# Technology-Stock-Data.xlsx / "itemloc" sheet edit
#
# The three columns D (max capacity), E (expiry date), F (fifo date) get
# rotated left by one position for the header row and every data row:
#   new D <- old E
#   new E <- old F
#   new F <- old D
#
# For the data rows, old E/F are always the literal text strings
# "2100-01-01" / "2000-01-01" (plain text, NOT real dates) and old D is a
# plain number. Writing a date-shaped literal straight into a cell via
# Range.Value/Value2 makes Excel auto-convert it into a real date serial
# (and also marks the cell with a quote-prefix style if forced to text),
# which would not match the original plain shared-string text cell. To
# avoid that, the date-looking strings are produced via a throwaway
# formula cell ( ="2100-01-01" ) whose computed value is pasted as
# values-only (PasteSpecial xlPasteValues) into the destination - this
# keeps the literal text, keeps the existing cell style (s="2"), and
# reuses the existing shared-string entries instead of minting new ones
# or new styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("itemloc")

# Scratch cells well outside the sheet's used range; cleared at the end.
$scratchExpiry = $ws.Cells.Item(200, 50)
$scratchFifo   = $ws.Cells.Item(200, 51)

# ---- Header row (row 1): rotate D1/E1/F1 ----
$ws.Cells.Item(1, 4).Value2 = "expiry date"
$ws.Cells.Item(1, 5).Value2 = "fifo date"
$ws.Cells.Item(1, 6).Value2 = "max capacity"

# ---- Data rows 2..63: rotate D/E/F ----
for ($r = 2; $r -le 63; $r++) {
    $oldD = $ws.Cells.Item($r, 4).Value2

    $scratchExpiry.Formula = "=""2100-01-01"""
    $scratchExpiry.Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4163)

    $scratchFifo.Formula = "=""2000-01-01"""
    $scratchFifo.Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4163)

    $ws.Cells.Item($r, 6).Value2 = $oldD
}

$scratchExpiry.Clear()
$scratchFifo.Clear()

# ---- Selection change: J8 -> J9 ----
[void]$ws.Range("J9").Select()
